# UPDATED usermodule -merged by vani
# Adds programId / batchId / userRoleProgramBatchStatus columns (U, V, W)
# to row 1 (headers) plus two data rows (2 and 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Create a throwaway named style so the new header font ("Aptos", 11pt,
# theme color 1) is added to the style table as plain direct formatting
# (fontId + applyFont="1") instead of being tracked as a permanent named
# cell style. Applying the named style assigns the new xf/font to the
# cells, and deleting the style afterwards collapses the xf back to the
# default xfId="0" while keeping the custom font - matching a normal
# "apply font" edit done directly in the Excel UI.
$headerStyleName = "TempHeaderFontStyle"
$headerStyle = $wb.Styles.Add($headerStyleName)
$headerStyle.Font.Name = "Aptos"
$headerStyle.Font.Size = 11

# --- New header cells ---
$ws.Range("U1").Value = "programId"
$ws.Range("V1").Value = "batchId"
$ws.Range("W1").Value = "userRoleProgramBatchStatus"
$ws.Range("U1:W1").Style = $headerStyleName

# Remove the temporary named style definition; the direct formatting
# (font) already applied to U1:W1 is preserved on the cells.
$wb.Styles($headerStyleName).Delete()

# --- New data rows ---
$ws.Range("U2").Value = 16224
$ws.Range("V2").Value = 8652
$ws.Range("W2").Value = "ACTIVE"

$ws.Range("U3").Value = 16210
$ws.Range("V3").Value = 8465
$ws.Range("W3").Value = "INACTIVE"

# --- Column widths for the newly introduced columns ---
$ws.Columns("U:W").AutoFit()

# --- Update selection / view to match the edited area ---
$ws.Range("U3").Select()
